# Applies the "Added more for testing" edit to Database.xlsx:
#  - Populates the RequestList sheet (sheet2) with a header row copied from
#    UserList plus one sample request row (banana / Chris / Moticska).
#  - Moves the active sheet/selection from UserList to RequestList.

$wb = $excel.ActiveWorkbook

$userList = $wb.Worksheets.Item("UserList")
$requestList = $wb.Worksheets.Item("RequestList")

# Header row, copied verbatim from UserList!A1:H1
$headers = @("username", "First Name", "Last Name", "DoB", "Card #", "Last Accessed", "Employee Status", "Password")
for ($col = 1; $col -le $headers.Length; $col++) {
    $requestList.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# New sample request row
$requestList.Range("A2").Value = "banana"
$requestList.Range("B2").Value = "Chris"
$requestList.Range("C2").Value = "Moticska"

# UserList's prior selection (B5) is cleared to a header-row range selection.
$userList.Range("A1:H1").Select() | Out-Null

# RequestList becomes the active/visible tab with its own selection.
$requestList.Activate()
$requestList.Range("F3").Select() | Out-Null
